$wb = $excel.ActiveWorkbook

# --- Rename existing sheets ---
$wb.Worksheets.Item(1).Name = "MultC"
$wb.Worksheets.Item(2).Name = "MultC_2"
$wb.Worksheets.Item(3).Name = "Quant"

# --- Add the new sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws4.Name = "MultAns"

# Write new unique strings in the exact order required so the shared
# string table matches the target (column A top-to-bottom first, then
# the two new comments in column C, row 6 before row 2).
$ws4.Range("A1").Value = "What is the best kind of bear?"
$ws4.Range("A2").Value = "Grizzly"
$ws4.Range("A3").Value = "Polar"
$ws4.Range("A4").Value = "Brown"
$ws4.Range("A6").Value = "Koala"
$ws4.Range("C6").Value = "It's not even a bear, dumbshit."
$ws4.Range("C2").Value = "Obviously it's not the grizzly."

# Remaining cells reuse already-existing shared strings.
$ws4.Range("A5").Value = "Black"

$ws4.Range("B1").Value = "Correct"
$ws4.Range("B2").Value = "N"
$ws4.Range("B3").Value = "Y"
$ws4.Range("B4").Value = "N"
$ws4.Range("B5").Value = "Y"
$ws4.Range("B6").Value = "N"

$ws4.Range("C1").Value = "Comment"
$ws4.Range("C3").Value = "You're brilliant!"
$ws4.Range("C4").Value = "You're stupid."
$ws4.Range("C5").Value = "So close. But not really that close."

# --- Column widths for the new sheet (closest achievable on Excel's column-width grid) ---
$ws4.Columns.Item(1).ColumnWidth = 34.5
$ws4.Columns.Item(2).ColumnWidth = 27.5
$ws4.Columns.Item(3).ColumnWidth = 38.833333333333336

# --- Sheet1 (MultC) selection update ---
$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Range("A1:C5").Select()

# --- Sheet2 (MultC_2) selection / active cell ---
$ws2 = $wb.Worksheets.Item(2)
[void]$ws2.Range("C6").Select()

# --- New sheet (MultAns) selection / active cell, and make it the active/visible tab ---
[void]$ws4.Range("C2").Select()
[void]$ws4.Activate()

Write-Host "done"
